$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.509.67'
$ws.Range("E2").Value = '  +2.04%  '
$ws.Range("D3").Value = '''3.558.57'
$ws.Range("E3").Value = '  +1.20%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''612.68'
$ws.Range("E5").Value = '  +6.01%  '
$ws.Range("D6").Value = '''173.13'
$ws.Range("E6").Value = '  +1.21%  '
$ws.Range("E7").Value = '  +1.52%  '
$ws.Range("D8").Value = '''3.554.03'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +4.05%  '
$ws.Range("D11").Value = '''7.24'
$ws.Range("E11").Value = '  +6.72%  '
$ws.Range("D12").Value = '''0.584'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").Value = '''46.53'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '''0.0000277'
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '''4.133.53'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").Value = '''8.36'
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").Value = '''615.24'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '''3.554.57'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '''70.563.44'
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("E20").Value = '  -2.09%  '
$ws.Range("D21").Value = '''17.34'
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").Value = '''9.38'
$ws.Range("E23").Value = '  -15.91%  '
$ws.Range("D24").Value = '''15.71'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").Value = '''96.73'
$ws.Range("E25").Value = '  -0.71%  '
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = '''33.47'
$ws.Range("E29").Value = '  +2.67%  '
$ws.Range("D30").Value = '''9.03'
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("D31").Value = '''8.49'
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("E33").Value = '  -1.52%  '
$ws.Range("D34").Value = '''6.96'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").Value = '''572.28'
$ws.Range("E35").Value = '  -9.73%  '
$ws.Range("E36").Value = '  +6.76%  '
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("D38").Value = '''10.81'
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").Value = '''57.26'
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("D40").Value = '''0.0471'
$ws.Range("E40").Value = '  +5.68%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +4.06%  '
$ws.Range("D43").Value = '''3.384.99'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("D45").Value = '''32.99'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '''2.96'
$ws.Range("E46").Value = '  +7.55%  '
$ws.Range("D47").Value = '''0.0₃0701'
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("E48").Value = '  +2.06%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '''133.86'
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("E51").Value = '  -0.01%  '

# Restore default (unstyled) formatting on the text-forced Price cells so
# the apostrophe-prefix trick does not leave a lingering quote-prefix style.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
